$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Add a note in G5 asking for the AP bond copy, styled as yellow text on a
# yellow fill (matches font/fill added to styles.xml in the target workbook).
$cell = $ws.Range("G5")
$cell.Value = "Send us the AP bond copy"
$cell.Font.Color = 65535
$cell.Interior.Color = 65535
